$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Dependencies"
$ws.Range("E1").Value = "Progress"

$ws.Range("E2").Value = 50
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("D6").Value = 4

$ws.Range("E3").Select()
